$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1621.8334
$ws.Range("J17").Value = 1661.3529
$ws.Range("L17").Value = 4984.0587
$ws.Range("N17").Value = -5320.0587
$ws.Range("H18").Value = 1698.9
$ws.Range("I18").Value = 1265.4445
$ws.Range("K18").Value = 1265.4445
$ws.Range("M18").Value = -981.4445000000001
$ws.Range("H19").Value = 3598.8
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 3598.8
$ws.Range("K19").Value = 0
$ws.Range("M19").Value = 3598.8
$ws.Range("N19").Value = -3948.8
$ws.Range("H31").Value = 453.77777
$ws.Range("J31").Value = 445.66666
$ws.Range("L31").Value = 1336.99998
$ws.Range("N31").Value = -1796.99998
$ws.Range("H32").Value = 3965.7778
$ws.Range("I32").Value = 3661.6
$ws.Range("J32").Value = 4346
$ws.Range("K32").Value = 3661.6
$ws.Range("L32").Value = 4346
$ws.Range("M32").Value = -3335.6
$ws.Range("N32").Value = -4998
$ws.Range("H39").Value = 269.86206
$ws.Range("J39").Value = 300.5
$ws.Range("L39").Value = 901.5
$ws.Range("N39").Value = -1493.5
$ws.Range("H55").Value = 1083.3334
$ws.Range("I55").Value = 425
$ws.Range("J55").Value = 1412.5
$ws.Range("K55").Value = 425
$ws.Range("L55").Value = 1412.5
$ws.Range("M55").Value = -211
$ws.Range("N55").Value = -1840.5
$ws.Range("H76").Value = 7396.25
$ws.Range("I76").Value = 4670
$ws.Range("J76").Value = 7785.7144
$ws.Range("K76").Value = 4670
$ws.Range("L76").Value = 7785.7144
$ws.Range("M76").Value = -4355
$ws.Range("N76").Value = -8415.714400000001
$ws.Range("H79").Value = 7396.25
$ws.Range("I79").Value = 4670
$ws.Range("J79").Value = 7785.7144
$ws.Range("K79").Value = 4670
$ws.Range("L79").Value = 7785.7144
$ws.Range("M79").Value = -3578
$ws.Range("N79").Value = -9969.714400000001
$ws.Range("H107").Value = 436.65
$ws.Range("I107").Value = 454.42105
$ws.Range("K107").Value = 454.42105
$ws.Range("M107").Value = 1465.57895
$ws.Range("H129").Value = 1156.6666
$ws.Range("I129").Value = 793.3333
$ws.Range("J129").Value = 2246.6667
$ws.Range("K129").Value = 2379.9999
$ws.Range("L129").Value = 6740.000100000001
$ws.Range("M129").Value = 2620.0001
$ws.Range("N129").Value = -16740.0001
$ws.Range("H132").Value = 1474.2727
$ws.Range("I132").Value = 1474.2727
$ws.Range("K132").Value = 4422.8181
$ws.Range("M132").Value = -1892.8181
$ws.Range("H135").Value = 2065.1738
$ws.Range("I135").Value = 2151.353
$ws.Range("J135").Value = 1821
$ws.Range("K135").Value = 19362.177
$ws.Range("L135").Value = 16389
$ws.Range("M135").Value = -16827.177
$ws.Range("N135").Value = -21459
$ws.Range("H138").Value = 2449.15
$ws.Range("I138").Value = 1149.9048
$ws.Range("J138").Value = 2911.5933
$ws.Range("K138").Value = 3449.7144
$ws.Range("L138").Value = 8734.7799
$ws.Range("M138").Value = 1690.2856
$ws.Range("N138").Value = -19014.7799

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26334696
$ws.Range("I32").Value = 38478316
$ws.Range("K32").Value = 38478316
$ws.Range("M32").Value = -38478029
$ws.Range("H45").Value = 62502850
$ws.Range("I45").Value = 125001440
$ws.Range("K45").Value = 125001440
$ws.Range("M45").Value = -125001063
$ws.Range("H123").Value = 47140
$ws.Range("J123").Value = 47140
$ws.Range("L123").Value = 47140
$ws.Range("N123").Value = -56940
$ws.Range("H127").Value = 119000
$ws.Range("J127").Value = 119000
$ws.Range("L127").Value = 119000
$ws.Range("N127").Value = -128920

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2921.6365
$ws.Range("I22").Value = 2013.8
$ws.Range("K22").Value = 2013.8
$ws.Range("M22").Value = -1840.8
$ws.Range("H80").Value = 2800.4443
$ws.Range("I80").Value = 7699
$ws.Range("J80").Value = 1400.8572
$ws.Range("K80").Value = 7699
$ws.Range("L80").Value = 1400.8572
$ws.Range("M80").Value = -6701
$ws.Range("N80").Value = -3396.8572
$ws.Range("H83").Value = 2800.4443
$ws.Range("I83").Value = 7699
$ws.Range("J83").Value = 1400.8572
$ws.Range("K83").Value = 38495
$ws.Range("L83").Value = 7004.286
$ws.Range("M83").Value = -33503
$ws.Range("N83").Value = -16988.286
$ws.Range("H86").Value = 2165.1428
$ws.Range("I86").Value = 2348.2778
$ws.Range("J86").Value = 1066.3334
$ws.Range("K86").Value = 2348.2778
$ws.Range("L86").Value = 1066.3334
$ws.Range("M86").Value = -1225.2778
$ws.Range("N86").Value = -3312.3334
$ws.Range("H89").Value = 2165.1428
$ws.Range("I89").Value = 2348.2778
$ws.Range("J89").Value = 1066.3334
$ws.Range("K89").Value = 11741.389
$ws.Range("L89").Value = 5331.666999999999
$ws.Range("M89").Value = -6125.388999999999
$ws.Range("N89").Value = -16563.667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 517
$ws.Range("I22").Value = 523.9
$ws.Range("K22").Value = 523.9
$ws.Range("M22").Value = -173.9
$ws.Range("H31").Value = 715634.5600000001
$ws.Range("I31").Value = 10543.632
$ws.Range("J31").Value = 1672543.8
$ws.Range("K31").Value = 10543.632
$ws.Range("L31").Value = 1672543.8
$ws.Range("M31").Value = -10248.632
$ws.Range("N31").Value = -1673133.8
$ws.Range("H34").Value = 715634.5600000001
$ws.Range("I34").Value = 10543.632
$ws.Range("J34").Value = 1672543.8
$ws.Range("K34").Value = 10543.632
$ws.Range("L34").Value = 1672543.8
$ws.Range("M34").Value = -10341.632
$ws.Range("N34").Value = -1672947.8
$ws.Range("H86").Value = 6060.5386
$ws.Range("I86").Value = 6198.3335
$ws.Range("J86").Value = 5942.4287
$ws.Range("K86").Value = 6198.3335
$ws.Range("L86").Value = 5942.4287
$ws.Range("M86").Value = -5075.3335
$ws.Range("N86").Value = -8188.4287
$ws.Range("H89").Value = 6060.5386
$ws.Range("I89").Value = 6198.3335
$ws.Range("J89").Value = 5942.4287
$ws.Range("K89").Value = 30991.6675
$ws.Range("L89").Value = 29712.1435
$ws.Range("M89").Value = -25375.6675
$ws.Range("N89").Value = -40944.14350000001
$ws.Range("H94").Value = 4737.8184
$ws.Range("I94").Value = 3111.6667
$ws.Range("J94").Value = 5347.625
$ws.Range("K94").Value = 3111.6667
$ws.Range("L94").Value = 5347.625
$ws.Range("M94").Value = -2660.6667
$ws.Range("N94").Value = -6249.625
$ws.Range("H122").Value = 1861.5
$ws.Range("I122").Value = 1861.5
$ws.Range("K122").Value = 5584.5
$ws.Range("M122").Value = -3134.5
$ws.Range("H132").Value = 3436.3845
$ws.Range("I132").Value = 2607.0908
$ws.Range("K132").Value = 7821.2724
$ws.Range("M132").Value = -5291.2724

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 656
$ws.Range("I11").Value = 656
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1968
$ws.Range("L11").Value = 0
$ws.Range("N11").Value = -1828

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 250
$ws.Range("J9").Value = 250
$ws.Range("L9").Value = 250
$ws.Range("N9").Value = -590

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 7250
$ws.Range("J58").Value = 8000
$ws.Range("L58").Value = 8000
$ws.Range("N58").Value = -8520
$ws.Range("H98").Value = 20000
$ws.Range("J98").Value = 20000
$ws.Range("L98").Value = 20000
$ws.Range("N98").Value = -25990
$ws.Range("H122").Value = 5218.8623
$ws.Range("I122").Value = 4800.0586
$ws.Range("J122").Value = 5812.1665
$ws.Range("K122").Value = 14400.1758
$ws.Range("L122").Value = 17436.4995
$ws.Range("M122").Value = -11950.1758
$ws.Range("N122").Value = -22336.4995
$ws.Range("H132").Value = 1445143.8
$ws.Range("I132").Value = 47502
$ws.Range("K132").Value = 142506
$ws.Range("M132").Value = -139976

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 45000
$ws.Range("J26").Value = 45000
$ws.Range("L26").Value = 45000
$ws.Range("N26").Value = -45586
$ws.Range("H126").Value = 2016.25
$ws.Range("I126").Value = 1711.375
$ws.Range("K126").Value = 5134.125
$ws.Range("M126").Value = -2664.125
$ws.Range("H132").Value = 917256.4399999999
$ws.Range("I132").Value = 10162.4
$ws.Range("K132").Value = 30487.2
$ws.Range("M132").Value = -27957.2
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120
